$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.063.37'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '1.557.45'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9996'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '287.46'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3858'
$ws.Range("E7").Value = '  +4.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3242'
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.46'
$ws.Range("E9").Value = '  -10.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.125'
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07362'
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9999'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.37'
$ws.Range("E13").Value = '  -5.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.704'
$ws.Range("E14").Value = '  -2.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.816'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '1.557.65'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001117'
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06607'
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.20'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.400'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9988'
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  -1.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.48'
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("D24").Value = '22.076.97'
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.333'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.542'
$ws.Range("E26").Value = '  -1.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.26'
$ws.Range("E27").Value = '  -0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.86'
$ws.Range("E28").Value = '  -2.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.860'
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").Value = '1.731.61'
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.78'
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.114'
$ws.Range("E32").Value = '  +6.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.866'
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.687'
$ws.Range("E34").Value = '  -14.44%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08188'
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.251'
$ws.Range("E36").Value = '  -4.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06231'
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02299'
$ws.Range("E38").Value = '  -4.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.218'
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2101'
$ws.Range("E40").Value = '  -4.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.221'
$ws.Range("E41").Value = '  -6.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.88'
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9989'
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5951'
$ws.Range("E44").Value = '  -2.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.49'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.720'
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5752'
$ws.Range("E47").Value = '  -3.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.930'
$ws.Range("E48").Value = '  -4.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.20'
$ws.Range("E49").Value = '  -3.71%  '
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06889'
$ws.Range("E51").Value = '  -3.78%  '
